# Weekly data refresh: insert the newest week's two Apio records
# (Primera / Segunda quality) at the top of the data block (row 641),
# pushing the previously-existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 641 (existing rows 641.. shift to 643..)
$ws.Rows.Item(641).Insert()
$ws.Rows.Item(641).Insert()

# --- Row 641: Apio, Primera, new week (Fecha serial 44706) ---
$ws.Cells.Item(641, 1).Value = 6
$ws.Cells.Item(641, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(641, 3).Value = "Metropolitana"
$ws.Cells.Item(641, 4).Value = 44706
$ws.Cells.Item(641, 5).Value = 13
$ws.Cells.Item(641, 6).Value = 100112017
$ws.Cells.Item(641, 7).Value = "Apio"
$ws.Cells.Item(641, 8).Value = "Americana (o)"
$ws.Cells.Item(641, 9).Value = "Primera"
$ws.Cells.Item(641, 10).Value = 2800
$ws.Cells.Item(641, 11).Value = 6000
$ws.Cells.Item(641, 12).Value = 7000
$ws.Cells.Item(641, 13).Value = 6571
$ws.Cells.Item(641, 14).Value = "$/docena de matas"
$ws.Cells.Item(641, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(641, 16).Value = 1095
$ws.Cells.Item(641, 17).Value = 6
$ws.Cells.Item(641, 18).Value = "Hortaliza"

# --- Row 642: Apio, Segunda, same new week ---
$ws.Cells.Item(642, 1).Value = 6
$ws.Cells.Item(642, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(642, 3).Value = "Metropolitana"
$ws.Cells.Item(642, 4).Value = 44706
$ws.Cells.Item(642, 5).Value = 13
$ws.Cells.Item(642, 6).Value = 100112017
$ws.Cells.Item(642, 7).Value = "Apio"
$ws.Cells.Item(642, 8).Value = "Americana (o)"
$ws.Cells.Item(642, 9).Value = "Segunda"
$ws.Cells.Item(642, 10).Value = 1602
$ws.Cells.Item(642, 11).Value = 400
$ws.Cells.Item(642, 12).Value = 5000
$ws.Cells.Item(642, 13).Value = 3277
$ws.Cells.Item(642, 14).Value = "$/docena de matas"
$ws.Cells.Item(642, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(642, 16).Value = 546
$ws.Cells.Item(642, 17).Value = 6
$ws.Cells.Item(642, 18).Value = "Hortaliza"
